$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (Column A and B for rows 2-4 get new test data)
$ws.Range("A2").Value = "TestingDemo"
$ws.Range("B2").Value = "Tes"

$ws.Range("A3").Value = "TestingDemoo"
$ws.Range("B3").Value = "Test"

$ws.Range("A4").Value = "TestingD"
$ws.Range("B4").Value = "TD"

# Update the selected cell in the sheet view
$ws.Range("B4").Select()
